# Auto-generated Excel COM-interop script that refreshes the cryptos worksheet
# with newly scraped coinranking.com figures (Coin name, Link, Price, Volume(1h)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / Link columns: a few rows got re-ranked, so the coin shown on that
#     row changes (e.g. row 44 was EnergySwap, now it is RenderToken). ---
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"

# --- Price column (D). Several new prices look like plain numbers (e.g. "308.35"),
#     which Excel would otherwise auto-convert from text to a numeric value. Force
#     each touched cell to Text format first, write the value, then restore the
#     cell style to "Normal" so no stray number formatting is left behind. ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.931.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.595.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.578"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.993.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.601.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.909"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "46.056.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "291.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.121"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.118.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "108.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.199"
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) column (E). These already contain surrounding spaces/% so Excel
#     keeps them as plain text automatically. ---
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("E23").Value = "  +15.08%  "
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  +2.63%  "
$ws.Range("E27").Value = "  +4.36%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("E32").Value = "  -3.13%  "
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("E35").Value = "  +3.73%  "
$ws.Range("E36").Value = "  -3.04%  "
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("E38").Value = "  -3.34%  "
$ws.Range("E39").Value = "  +3.88%  "
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("E44").Value = "  -4.22%  "
$ws.Range("E45").Value = "  +5.13%  "
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("E48").Value = "  +5.08%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  -0.62%  "
